$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new column F, matching style of existing header cell E1
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F1").Value = "time_taken"

# Per-row time_taken values (microsecond-precision timestamps)
$times = @(
    "2021-10-05 10:52:15.181017",
    "2021-10-05 10:52:15.181028",
    "2021-10-05 10:52:15.181032",
    "2021-10-05 10:52:15.181035",
    "2021-10-05 10:52:15.181039",
    "2021-10-05 10:52:15.181042",
    "2021-10-05 10:52:15.181045",
    "2021-10-05 10:52:15.181048",
    "2021-10-05 10:52:15.181051",
    "2021-10-05 10:52:15.181054"
)

for ($i = 0; $i -lt $times.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $times[$i]
}
